$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original wide-format fold-change values before they get
# overwritten while reshaping the sheet into long format.
$dusp11Value = 0.6659568311637312
$ifnbValue   = 2.215936028923374
$mx1Value    = 1.170718510326486

# Propagate the header/number style (s="1": bold, bordered, centered) from
# A2 down to A3:A4 before we touch anything else.
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New column headers
$ws.Range("B1").Value = "Condition"
$ws.Range("C1").Value = "Fold Change"

# Drop the now-unused D column (data moved into long format below)
$ws.Range("D1:D2").Clear()

# Row 2: dusp11
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "dusp11 foldchange"
$ws.Cells.Item(2, 3).Value = $dusp11Value

# Row 3: ifnb
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "ifnb foldchange"
$ws.Cells.Item(3, 3).Value = $ifnbValue

# Row 4: mx1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "mx1 foldchange"
$ws.Cells.Item(4, 3).Value = $mx1Value
